# Apply updated cryptocurrency price/volume figures pulled by the
# scheduled GitHub Actions refresh job. Row 45/46 additionally swap
# (Stellar and VeChain traded ranking positions).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ column letter = new value }
$rowUpdates = [ordered]@{
    2 = @{ 'D'="52.585.85"; 'E'="  -13.15%  " }
    3 = @{ 'D'="2.309.34"; 'E'="  -20.36%  " }
    4 = @{ 'E'="  +0.16%  " }
    5 = @{ 'D'="447.84"; 'E'="  -14.90%  " }
    6 = @{ 'D'="121.47"; 'E'="  -14.91%  " }
    7 = @{ 'D'="0.998"; 'E'="  -0.11%  " }
    8 = @{ 'D'="0.468"; 'E'="  -14.65%  " }
    9 = @{ 'D'="2.312.48"; 'E'="  -20.45%  " }
    10 = @{ 'D'="5.30"; 'E'="  -11.10%  " }
    11 = @{ 'D'="0.0869"; 'E'="  -18.68%  " }
    12 = @{ 'D'="0.302"; 'E'="  -15.72%  " }
    13 = @{ 'E'="  -6.02%  " }
    14 = @{ 'D'="52.573.30"; 'E'="  -13.17%  " }
    15 = @{ 'D'="18.82"; 'E'="  -16.47%  " }
    16 = @{ 'D'="0.0000118"; 'E'="  -16.10%  " }
    17 = @{ 'D'="2.322.00"; 'E'="  -20.13%  " }
    18 = @{ 'D'="3.96"; 'E'="  -20.32%  " }
    19 = @{ 'D'="297.79"; 'E'="  -15.02%  " }
    20 = @{ 'D'="8.91"; 'E'="  -23.15%  " }
    21 = @{ 'E'="  -0.16%  " }
    22 = @{ 'E'="  -1.38%  " }
    23 = @{ 'D'="5.14"; 'E'="  -21.14%  " }
    24 = @{ 'D'="53.56"; 'E'="  -17.08%  " }
    25 = @{ 'D'="0.365"; 'E'="  -19.18%  " }
    26 = @{ 'D'="0.144"; 'E'="  -18.86%  " }
    27 = @{ 'D'="6.92"; 'E'="  -11.21%  " }
    28 = @{ 'D'="0.998"; 'E'="  -0.14%  " }
    29 = @{ 'D'="0.0₃0661"; 'E'="  -22.05%  " }
    30 = @{ 'D'="140.92"; 'E'="  -6.79%  " }
    31 = @{ 'D'="16.83"; 'E'="  -14.01%  " }
    32 = @{ 'E'="  -19.71%  " }
    33 = @{ 'D'="4.72"; 'E'="  -15.09%  " }
    34 = @{ 'E'="  -17.59%  " }
    35 = @{ 'D'="3.40"; 'E'="  -21.11%  " }
    36 = @{ 'D'="0.994"; 'E'="  -0.32%  " }
    37 = @{ 'E'="  -16.93%  " }
    38 = @{ 'D'="31.81"; 'E'="  -15.54%  " }
    39 = @{ 'D'="10.15"; 'E'="  -1.66%  " }
    40 = @{ 'E'="  -13.48%  " }
    41 = @{ 'E'="  -13.03%  " }
    42 = @{ 'E'="  -16.12%  " }
    43 = @{ 'D'="1.908.27"; 'E'="  -16.65%  " }
    44 = @{ 'E'="  -20.16%  " }
    45 = @{ 'B'="Stellar"; 'C'="https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; 'D'="0.0820"; 'E'="  -10.63%  " }
    46 = @{ 'B'="VeChain"; 'C'="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; 'D'="0.0204"; 'E'="  -13.64%  " }
    47 = @{ 'D'="4.15"; 'E'="  -15.94%  " }
    48 = @{ 'D'="15.53"; 'E'="  -23.82%  " }
    49 = @{ 'E'="  -5.19%  " }
    50 = @{ 'D'="4.45"; 'E'="  -12.78%  " }
    51 = @{ 'D'="14.96"; 'E'="  -18.04%  " }
}

foreach ($row in $rowUpdates.Keys) {
    foreach ($col in $rowUpdates[$row].Keys) {
        $cellRef = "$col$row"
        $newValue = $rowUpdates[$row][$col]
        if ($col -eq "D") {
            # Column D holds price text (e.g. "52.585.85", "0.998") that must
            # stay text -- force the Text format so Excel does not reinterpret
            # numeric-looking values (and strip meaningful trailing zeros).
            $ws.Range($cellRef).NumberFormat = "@"
        }
        $ws.Range($cellRef).Value = $newValue
    }
}
